$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K is a copy of column J (same "status" header + PASS/FAIL values,
# same per-row formatting) — mirrors how the J column itself was built.
$ws.Range("J1:J6").Copy() | Out-Null
$ws.Range("K1").PasteSpecial() | Out-Null

# Re-apply the header fill explicitly so K1 keeps the same highlighted
# look as J1 (and the rest of the "status" header row).
$ws.Range("K1").Interior.Color = $ws.Range("J1").Interior.Color

# Match column K's width to the other "status" columns (C:J).
$ws.Columns.Item(11).ColumnWidth = 5.52
